$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25 (pushes rows 25-34 down to 26-35)
$ws.Rows("25").Insert()

# Copy formatting (styles/merges layout) from the row below (now row 26, former row 25)
$ws.Range("A26:Q26").Copy()
$ws.Range("A25:Q25").PasteSpecial(-4122)

# Recreate the merged cells for the new row 25 to match the other data rows
$ws.Range("A25:B25").Merge()
$ws.Range("C25:G25").Merge()
$ws.Range("H25:K25").Merge()
$ws.Range("L25:M25").Merge()
$ws.Range("N25:O25").Merge()

# Fill in the new item's data
# (values that look numeric are prefixed with an apostrophe so they are stored
# as text, matching how this workbook stores these columns as shared strings)
$ws.Range("A25").Value = 19
$ws.Range("C25").Value = "بيرسول حشرات طائره"
$ws.Range("H25").Value = "9:0"
$ws.Range("L25").Value = "'0"
$ws.Range("N25").Value = "'50.00"
$ws.Range("P25").Value = "'50.0000"
$ws.Range("Q25").Value = "'1:0"

# Renumber the sequential "م" column for the rows that shifted down
$ws.Range("A26").Value = 20
$ws.Range("A27").Value = 21
$ws.Range("A28").Value = 22
$ws.Range("A29").Value = 23
$ws.Range("A30").Value = 24
$ws.Range("A31").Value = 25
$ws.Range("A32").Value = 26
$ws.Range("A33").Value = 27

# Update the grand total (row 34, column P) to include the new item's sell price
$ws.Range("P34").Value = 1485.0650000000001

# Update the generated timestamp in the footer (row 35)
$ws.Range("A35").Value = "Friday, 20 June, 2025 7:04 PM"
